$d = $word.ActiveDocument

# --- 1) Rewrite the ATTRS processing-instruction paragraph (adds ",featured" with proofing marks) ---
$attrsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HEDProcessinginstruction"/></w:pPr><w:r><w:t xml:space="preserve">ATTRS=id: selectors; data-tags: </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>typeset</w:t></w:r><w:r><w:t>,featured</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>;</w:t></w:r></w:p>
'@
$d.Paragraphs(2).Range.InsertXML($attrsXml)

# --- 2) Rewrite the chapter/section title paragraph ---
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HEDChapterSectiontitle"/></w:pPr><w:r><w:t>Use the “Limit these changes” Menu to Adjust</w:t></w:r><w:r><w:t xml:space="preserve"> the </w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t xml:space="preserve">esign of </w:t></w:r><w:r><w:t>O</w:t></w:r><w:r><w:t xml:space="preserve">nly </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t xml:space="preserve">ertain </w:t></w:r><w:r><w:t>P</w:t></w:r><w:r><w:t>aragraphs</w:t></w:r><w:r><w:t xml:space="preserve"> or Elements</w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($titleXml)

# --- 3) Append three new paragraphs after the trailing "HEDBOXBoxend" paragraph ---
$videoXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HEDPlaintextparagraph"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Here’s a video that walks you through the process of using the “Limit these changes” menu to apply a design to a single paragraph:</w:t></w:r></w:p>
'@
$newPara1 = $d.Paragraphs.Add()
$newPara1.Range.InsertXML($videoXml)

$iframeXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HEDPlaintextparagraph"/></w:pPr><w:r><w:t xml:space="preserve">&lt;iframe width="560" height="315" </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="https://www.youtube.com/embed/HrpE181HFd8" frameborder="0" allow="accelerometer; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>autoplay</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">; clipboard-write; encrypted-media; gyroscope; picture-in-picture" </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>allowfullscreen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;&lt;/iframe&gt;</w:t></w:r></w:p>
'@
$newPara2 = $d.Paragraphs.Add()
$newPara2.Range.InsertXML($iframeXml)

$htmlXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HEDProcessinginstruction"/></w:pPr><w:r><w:t>HTML=true</w:t></w:r></w:p>
'@
$newPara3 = $d.Paragraphs.Add()
$newPara3.Range.InsertXML($htmlXml)

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"
